# Applies the OOXML diff: splits the declaration paragraph's single run
# into three runs so that "ação {{ tipo_acao }}" becomes
# "{{ tipo_acao.upper() }}" (rg/tipoacao variable change).
$d = $word.ActiveDocument

# Locate the paragraph that still holds the un-split "tipo_acao" placeholder
# (robust to the paragraph's ordinal position in the document).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*tipo_acao*") {
        $target = $candidate
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the declaration paragraph containing 'tipo_acao'."
}

# Replace the paragraph's run content (but not its end-of-paragraph mark,
# so w:pPr / w14:paraId and friends on the <w:p> survive untouched) with the
# new three-run markup via InsertXML, which lets us control run boundaries
# exactly - a plain Range.Text/Find.Execute replace always re-collapses the
# paragraph back down to a single run.
$pStart = $target.Range.Start
$pEnd = $target.Range.End
$contentRange = $d.Range($pStart, $pEnd - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:xml="http://www.w3.org/XML/1998/namespace"><w:body><w:p><w:r w:rsidRPr="00963658"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/><w:t>Eu, {{ nome_assistido }}, CPF nº {{ num_cpf }}, RG nº {{ num_rg }}, declaro, para fins de direito, sob as penas da lei, que as informações prestadas são fiéis à verdade e condizentes com a realidade dos fatos à época. Declaro ainda que tenho ciência do teor da petição inicial para a {{ tipo_acao</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.upper()</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> }}, do Caso {{ num_caso }}, concordando com o que ali foi exposto e, consequentemente, com a distribuição da ação. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$contentRange.InsertXML($xml)

Write-Output ("Updated paragraph text: [" + $target.Range.Text + "]")
